$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the header cell G1 from "Template_name" to "Template_name_strategy"
$ws.Range("G1").Value = "Template_name_strategy"

# Move the active selection to G1 (from G5)
$ws.Range("G1").Select()
